$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update column C ("Förändrad") from 2023-09-06 (45175) to 2023-09-14 (45183)
# for all data rows (rows 2 through 70).
$newDate = (Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0 -Millisecond 0).AddDays(45183)

for ($row = 2; $row -le 70; $row++) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
